# Update the NATMI Fgf17-Fgfr2 LR-pairs table with the newly re-run TPM
# values. The sending/target cluster grid gained a 4th target cluster
# ("Resolving-Mac"), so rows 2-9 are refreshed in place and rows 10-13 are
# appended, mirroring the source script's regenerated output table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('ECs', 'Fgf17', 'Fgfr2', 'ECs', 1, 0.3333333333333333, 0.09982999999999999, 0.29949, 0.06953924013247029, 0.06953924013247029, 3, 1, 0.7144740000000001, 2.143422, 0.138796410342318, 0.138796410342318, 0.07132593942, 0.64193345478, 0.009651796908319335, 0.009651796908319335),
    @('ECs', 'Fgf17', 'Fgfr2', 'FAPs', 1, 0.3333333333333333, 0.09982999999999999, 0.29949, 0.06953924013247029, 0.06953924013247029, 3, 1, 4.140873, 12.422619, 0.8044215857867821, 0.8044215857867821, 0.4133833515899999, 3.72045016431, 0.05593886582176959, 0.05593886582176959),
    @('ECs', 'Fgf17', 'Fgfr2', 'MuSCs', 1, 0.3333333333333333, 0.09982999999999999, 0.29949, 0.06953924013247029, 0.06953924013247029, 3, 1, 0.2847646666666667, 0.8542940000000001, 0.05531945672713084, 0.05531945672713083, 0.02842805667333334, 0.25585251006, 0.003846872985345751, 0.00384687298534575),
    @('ECs', 'Fgf17', 'Fgfr2', 'Resolving-Mac', 1, 0.3333333333333333, 0.09982999999999999, 0.29949, 0.06953924013247029, 0.06953924013247029, 2, 0.6666666666666666, 0.007528666666666667, 0.022586, 0.00146254714376898, 0.00146254714376898, 0.0007515867933333333, 0.00676428114, 0.0001017044170356096, 0.0001017044170356096),
    @('FAPs', 'Fgf17', 'Fgfr2', 'ECs', 3, 1, 0.6306993333333334, 1.892098, 0.4393303855760352, 0.4393303855760352, 3, 1, 0.7144740000000001, 2.143422, 0.138796410342318, 0.138796410342318, 0.4506182754840001, 4.055564479356001, 0.06097748047226018, 0.06097748047226018),
    @('FAPs', 'Fgf17', 'Fgfr2', 'FAPs', 3, 1, 0.6306993333333334, 1.892098, 0.4393303855760352, 0.4393303855760352, 3, 1, 4.140873, 12.422619, 0.8044215857867821, 0.8044215857867821, 2.611645840518, 23.504812564662, 0.3534068454493927, 0.3534068454493927),
    @('FAPs', 'Fgf17', 'Fgfr2', 'MuSCs', 3, 1, 0.6306993333333334, 1.892098, 0.4393303855760352, 0.4393303855760352, 3, 1, 0.2847646666666667, 0.8542940000000001, 0.05531945672713084, 0.05531945672713083, 0.1796008854235556, 1.616407968812001, 0.02430351825378719, 0.02430351825378718),
    @('FAPs', 'Fgf17', 'Fgfr2', 'Resolving-Mac', 3, 1, 0.6306993333333334, 1.892098, 0.4393303855760352, 0.4393303855760352, 2, 0.6666666666666666, 0.007528666666666667, 0.022586, 0.00146254714376898, 0.00146254714376898, 0.004748325047555556, 0.04273492542800001, 0.000642541400595155, 0.0006425414005951549),
    @('MuSCs', 'Fgf17', 'Fgfr2', 'ECs', 3, 1, 0.705063, 2.115189, 0.4911303742914945, 0.4911303742914945, 3, 1, 0.7144740000000001, 2.143422, 0.138796410342318, 0.138796410342318, 0.503749181862, 4.533742636758, 0.06816713296173851, 0.06816713296173851),
    @('MuSCs', 'Fgf17', 'Fgfr2', 'FAPs', 3, 1, 0.705063, 2.115189, 0.4911303742914945, 0.4911303742914945, 3, 1, 4.140873, 12.422619, 0.8044215857867821, 0.8044215857867821, 2.919576339999, 26.276187059991, 0.3950758745156199, 0.3950758745156199),
    @('MuSCs', 'Fgf17', 'Fgfr2', 'MuSCs', 3, 1, 0.705063, 2.115189, 0.4911303742914945, 0.4911303742914945, 3, 1, 0.2847646666666667, 0.8542940000000001, 0.05531945672713084, 0.05531945672713083, 0.200777030174, 1.806993271566, 0.02716906548799791, 0.0271690654879979),
    @('MuSCs', 'Fgf17', 'Fgfr2', 'Resolving-Mac', 3, 1, 0.705063, 2.115189, 0.4911303742914945, 0.4911303742914945, 2, 0.6666666666666666, 0.007528666666666667, 0.022586, 0.00146254714376898, 0.00146254714376898, 0.005308184306, 0.047773658754, 0.0007183013261382155, 0.0007183013261382154)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
